$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 observation (added first so the new shared string lands at the same
# index the canonical workbook uses)
$ws.Range("H18").Value = "Pendiente completar para pasar a montaje"

# Área: Ciencias Naturales
$ws.Range("C1").Value = "Ciencias Naturales"

# Row 16 observation
$ws.Range("F16").Value = "En manuscritos"

# Grado: 7
$ws.Range("C2").Value = 7

# Update the active selection to F16 (matches the diff's sheetView/selection change)
$ws.Range("F16").Select()
